# Applies the 21-11-2023 20:45 script update to the England League One
# 2023-2024 sheet:
#   1) A permutation of the match-data columns (F:V) across a set of
#      existing rows (row index / country / tournament / season / date in
#      A:E are untouched - only which match's odds/teams sit in each row
#      changes).
#   2) One new match appended as row 200 (Leyton Orient 0-1 Lincoln).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target row -> source row: the F:V payload that ends up in the target row
# is whatever currently (pre-edit) sits in the source row.
$map = @{
  93  = 96
  96  = 93
  101 = 102
  102 = 101
  103 = 107
  104 = 113
  105 = 112
  106 = 111
  107 = 110
  108 = 109
  109 = 114
  110 = 108
  111 = 106
  112 = 105
  113 = 104
  114 = 103
  115 = 117
  116 = 119
  117 = 115
  118 = 116
  119 = 118
  120 = 124
  121 = 120
  122 = 125
  124 = 122
  125 = 121
  127 = 136
  128 = 127
  129 = 131
  130 = 128
  131 = 129
  132 = 130
  133 = 135
  134 = 132
  135 = 133
  136 = 134
  142 = 146
  143 = 142
  146 = 143
  160 = 164
  161 = 160
  162 = 161
  163 = 162
  164 = 163
  166 = 167
  167 = 168
  168 = 166
  169 = 171
  171 = 169
}

# Snapshot every involved row's F:V values BEFORE writing anything, since
# the mapping above is made up of several disjoint permutation cycles
# (not just 2-element swaps) - writing in place while reading would clobber
# data still needed for a later row in the same cycle.
$snapshot = @{}
foreach ($row in $map.Keys) {
  $snapshot[$row] = $ws.Range("F$row`:V$row").Value()
}

foreach ($row in $map.Keys) {
  $source = $map[$row]
  $ws.Range("F$row`:V$row").Value = $snapshot[$source]
}

# Append the new match as row 200. Clone the formatting (bold/bordered
# index style + date-time number format on E) from row 199, the previous
# last row, before filling in the values.
$ws.Range("A199:E199").Copy() | Out-Null
$ws.Range("A200:E200").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(200, 1).Value = 199
$ws.Cells.Item(200, 2).Value = "england"
$ws.Cells.Item(200, 3).Value = "league-one"
$ws.Cells.Item(200, 4).Value = "2023-2024"
$ws.Cells.Item(200, 5).Value = 45251.86458333334
$ws.Cells.Item(200, 6).Value = "Leyton Orient"
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = "Lincoln"
$ws.Cells.Item(200, 9).Value = 1
$ws.Cells.Item(200, 10).Value = 2.49
$ws.Cells.Item(200, 11).Value = "30/09/2023 23:42"
$ws.Cells.Item(200, 12).Value = 1.9
$ws.Cells.Item(200, 13).Value = "21/11/2023 20:30"
$ws.Cells.Item(200, 14).Value = 3.16
$ws.Cells.Item(200, 15).Value = "30/09/2023 23:42"
$ws.Cells.Item(200, 16).Value = 3.38
$ws.Cells.Item(200, 17).Value = "21/11/2023 20:30"
$ws.Cells.Item(200, 18).Value = 2.92
$ws.Cells.Item(200, 19).Value = "30/09/2023 23:42"
$ws.Cells.Item(200, 20).Value = 4.68
$ws.Cells.Item(200, 21).Value = "21/11/2023 20:30"
$ws.Cells.Item(200, 22).Value = "https://www.betexplorer.com/football/england/league-one/leyton-orient-lincoln-city/t6fmrCI2/"
